# Update "想去人数" (number of people wanting to go) counts on both the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet, which
# duplicates the same rows with one extra row inserted near the bottom.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value = 218
$wsExhibit.Range("F8").Value = 14559
$wsExhibit.Range("F9").Value = 163
$wsExhibit.Range("F11").Value = 5813
$wsExhibit.Range("F15").Value = 65
$wsExhibit.Range("F18").Value = 83
$wsExhibit.Range("F23").Value = 10602

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 218
$wsAll.Range("F9").Value = 14559
$wsAll.Range("F10").Value = 163
$wsAll.Range("F12").Value = 5813
$wsAll.Range("F16").Value = 65
$wsAll.Range("F19").Value = 83
$wsAll.Range("F25").Value = 10602
